$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text (not auto-number) formatting is preserved for D/E columns as inline strings

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.888.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.696.16'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.01'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.19'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.52%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.125'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.99'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.31%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '30.70'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000208'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +10.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.182.32'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.704.32'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.693.29'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.74'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.70'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.70%  '
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '359.83'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.16'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.90'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.82%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +12.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.67'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.67'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.172'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.20%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.49'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.80%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.20'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.98%  '
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '544.82'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.81'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.69'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.47'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.435'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.83'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '164.14'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.79%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '42.67'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '168.67'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.21'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0617'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.67'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.29'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.658'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.86'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +7.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0993'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.02%  '
